$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.077.44'
$ws.Range('E2').Value = '  +1.99%  '

$ws.Range('D3').Value = '3.621.98'
$ws.Range('E3').Value = '  +3.07%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '605.13'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.43%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '201.91'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.53%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.631'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.24%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.04%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.218'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +7.88%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.649'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.85%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.42'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.54%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000305'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.47%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.62'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.62%  '

$ws.Range('D14').Value = '4.194.08'
$ws.Range('E14').Value = '  +3.00%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '683.20'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +14.61%  '

$ws.Range('D16').Value = '71.045.28'
$ws.Range('E16').Value = '  +1.73%  '

$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.633.31'
$ws.Range('E17').Value = '  +3.00%  '

$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.84'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.25%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '19.15'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.03%  '

$ws.Range('E20').Value = '  +0.43%  '

$ws.Range('E21').Value = '  +1.88%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '19.14'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.58%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.44'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.38%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '106.90'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.02%  '

$ws.Range('E25').Value = '  +0.49%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.04'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.68%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.60'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.06%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.99'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.93%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.15'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.92%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.27'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.52%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.56'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.25%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.32'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.64%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.27'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.06%  '

$ws.Range('E34').Value = '  +1.17%  '

$ws.Range('E35').Value = '  +0.71%  '

$ws.Range('D36').Value = '0.0₃0864'
$ws.Range('E36').Value = '  +5.09%  '

$ws.Range('D37').Value = '3.905.40'
$ws.Range('E37').Value = '  +4.39%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '530.11'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.55%  '

$ws.Range('E39').Value = '  -0.05%  '

$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.06'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.61%  '

$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.61'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.39%  '

$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '36.98'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.32%  '

$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.391'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.28%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.141'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.27%  '

$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0468'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.43%  '

$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.06'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +9.07%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.46'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +6.67%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.142'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.58%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.66'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.64%  '

$ws.Range('E50').Value = '  -0.33%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.000245'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.68%  '
